$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# New row 36: GoodSplit / maxLen:2 result for TeacherRL
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = 602
$ws.Cells.Item(36, 3).Value = 602
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = 2631
$ws.Cells.Item(36, 6).Value = "../data/tests/sequences/Moore_R10_PDS.fsm"
$ws.Cells.Item(36, 7).Value = "GoodSplit"
$ws.Cells.Item(36, 8).Value = "maxLen:2"
$ws.Cells.Item(36, 9).Value = "TeacherRL"

# New row 37: GoodSplit / maxLen:2 + EQtoStop result for TeacherRL
$ws.Cells.Item(37, 1).Value = 1
$ws.Cells.Item(37, 2).Value = 199
$ws.Cells.Item(37, 3).Value = 199
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(37, 5).Value = 752
$ws.Cells.Item(37, 6).Value = "../data/tests/sequences/Moore_R10_PDS.fsm"
$ws.Cells.Item(37, 7).Value = "GoodSplit"
$ws.Cells.Item(37, 8).Value = "maxLen:2 + EQtoStop"
$ws.Cells.Item(37, 9).Value = "TeacherRL"

# Move the view / selection to the newly added rows
$ws.Range("A36:I37").Select()
